# Weekly update: add the newest Femacal de La Calera - Papaya price record.
# A new row is inserted at row 40 (the sheet's records are kept newest-first
# among the most-recent entries), pushing the previous rows 40:73 down to
# 41:74 and growing the used range from A1:T73 to A1:T74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40, shifting rows 40:73 -> 41:74.
$ws.Rows.Item(40).Insert()

# Fill in the new weekly record.
$ws.Range('A40').Value = 3
$ws.Range('B40').Value = 'Femacal de La Calera'
$ws.Range('C40').Value = 'Coquimbo'
$ws.Range('D40').Value = 45126
$ws.Range('E40').Value = 5
$ws.Range('F40').Value = 'Fruta'
$ws.Range('G40').Value = 100108
$ws.Range('H40').Value = 'Tropicales y subtropicales'
$ws.Range('I40').Value = 100108004
$ws.Range('J40').Value = 'Papaya'
$ws.Range('K40').Value = 'Cultivar IV Región'
$ws.Range('L40').Value = 'Primera'
$ws.Range('M40').Value = 56
$ws.Range('N40').Value = 20000
$ws.Range('O40').Value = 20000
$ws.Range('P40').Value = 20000
$ws.Range('Q40').Value = '$/bandeja 10 kilos'
$ws.Range('R40').Value = 'Provincia del Elquí'
$ws.Range('S40').Value = 2000
$ws.Range('T40').Value = 10
